# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 3 everywhere corresponds to the c0fce3fe... file which was
# "Ready for handoff"; the handback transform for it has now failed.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Widen the "Error Detail" column (P) on both locale sheets to fit the new
# error messages (stored column width ends up 40 characters).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664

# Populate the Error Detail column for row 3 (the c0fce3fe... file) on each
# locale sheet with the handback-transform error message.
$wsZhCn.Range("P3").Value = "Handback file name: 35zu4edv.msd is different with handoff file name: c0fce3fe-50d6-4bea-8068-d9a62471ce85.76e05161683fccd48fc3e70678f93d5be0121865.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: 35zu4edv.msd is different with handoff file name: c0fce3fe-50d6-4bea-8068-d9a62471ce85.76e05161683fccd48fc3e70678f93d5be0121865.de-de."
